# "Nuevo formato 15 jun 2021" - updated grade distribution / pass-rate figures
# across the five summary sheets (grades were reshuffled between buckets,
# which also shifts the dependent Repro/Por_Repro and Apro*/Por_Apro* totals).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Totales Plantel 1P")
$ws.Range("F5").Value = 5
$ws.Range("G5").Value = 7
$ws.Range("I5").Value = 21
$ws.Range("J5").Value = 18
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 10
$ws.Range("G6").Value = 17
$ws.Range("H6").Value = 20
$ws.Range("I6").Value = 13
$ws.Range("K6").Value = 81
$ws.Range("L6").Value = 95
$ws.Range("M6").Value = 53.98
$ws.Range("G7").Value = 5
$ws.Range("H7").Value = 12

$ws = $wb.Worksheets.Item("Totales Plantel 2P")
$ws.Range("I2").Value = 21
$ws.Range("J2").Value = 23
$ws.Range("F5").Value = 8
$ws.Range("H5").Value = 19
$ws.Range("J5").Value = 20
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 19
$ws.Range("G6").Value = 13
$ws.Range("H6").Value = 24
$ws.Range("I6").Value = 15
$ws.Range("J6").Value = 22
$ws.Range("K6").Value = 80
$ws.Range("L6").Value = 96
$ws.Range("M6").Value = 54.55
$ws.Range("G7").Value = 8
$ws.Range("H7").Value = 9

$ws = $wb.Worksheets.Item("Totales Plantel Final")
$ws.Range("F2").Value = 2
$ws.Range("H2").Value = 15
$ws.Range("I2").Value = 16
$ws.Range("J2").Value = 25
$ws.Range("K2").Value = 136
$ws.Range("L2").Value = 79
$ws.Range("M2").Value = 36.74
$ws.Range("G5").Value = 7
$ws.Range("H5").Value = 5
$ws.Range("J5").Value = 19
$ws.Range("K5").Value = 94
$ws.Range("L5").Value = 58
$ws.Range("M5").Value = 38.16
$ws.Range("F6").Value = 4
$ws.Range("G6").Value = 9
$ws.Range("I6").Value = 22
$ws.Range("J6").Value = 21
$ws.Range("K6").Value = 106
$ws.Range("L6").Value = 70
$ws.Range("M6").Value = 39.77
$ws.Range("G7").Value = 3
$ws.Range("J7").Value = 11

$ws = $wb.Worksheets.Item("Reprobados por Grupo")
$ws.Range("F3").Value = 1
$ws.Range("I3").Value = 8
$ws.Range("J3").Value = 28
$ws.Range("K3").Value = 13
$ws.Range("L3").Value = 31.71
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 3
$ws.Range("G7").Value = 2
$ws.Range("H7").Value = 3
$ws.Range("I7").Value = 3
$ws.Range("F22").Value = 0
$ws.Range("G22").Value = 0
$ws.Range("I22").Value = 1
$ws.Range("J22").Value = 17
$ws.Range("K22").Value = 4
$ws.Range("L22").Value = 19.05
$ws.Range("E26").Value = 0
$ws.Range("G26").Value = 0
$ws.Range("I26").Value = 2
$ws.Range("J26").Value = 27
$ws.Range("K26").Value = 6
$ws.Range("L26").Value = 18.18
$ws.Range("E27").Value = 0
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 7
$ws.Range("I27").Value = 4
$ws.Range("J27").Value = 22
$ws.Range("K27").Value = 17
$ws.Range("L27").Value = 43.59
$ws.Range("H28").Value = 3
$ws.Range("I28").Value = 8
$ws.Range("F29").Value = 0
$ws.Range("G29").Value = 5
$ws.Range("H29").Value = 5
$ws.Range("I29").Value = 2
$ws.Range("J29").Value = 24
$ws.Range("K29").Value = 12
$ws.Range("L29").Value = 33.33
$ws.Range("F32").Value = 1
$ws.Range("I32").Value = 5

$ws = $wb.Worksheets.Item("Totales Grupos")
$ws.Range("G3").Value = 28
$ws.Range("H3").Value = 68.29
$ws.Range("G22").Value = 17
$ws.Range("H22").Value = 80.95
$ws.Range("E26").Value = 21
$ws.Range("F26").Value = 63.64
$ws.Range("G26").Value = 27
$ws.Range("H26").Value = 81.82
$ws.Range("C27").Value = 20
$ws.Range("D27").Value = 51.28
$ws.Range("G27").Value = 22
$ws.Range("H27").Value = 56.41
$ws.Range("C29").Value = 21
$ws.Range("D29").Value = 58.33
$ws.Range("G29").Value = 24
$ws.Range("H29").Value = 66.67
